# Refresh cryptocurrency prices and 1h volume-change percentages.
# (GitHub Actions scheduled data update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: Price column (D) values are plain-looking decimal numbers that
# Excel would otherwise auto-convert to numeric type (dropping trailing zeros),
# so we force Text format before assignment, then restore the default style
# so no stray per-cell formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.730.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.360.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.30%  "
$ws.Range("E7").Value = "  -2.20%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0927"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.18%  "
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.714.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.550.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.687.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.22%  "
$ws.Range("E24").Value = "  -3.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0894"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("E34").Value = "  -8.40%  "
$ws.Range("E35").Value = "  +17.03%  "
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0365"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("E39").Value = "  -5.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.240"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  -5.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "112.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.62%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "
